$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66; this shifts existing rows 66-74 down to 67-75
$ws.Rows("66").Insert()

# Populate the newly inserted row 66 with the new price record
$ws.Range("A66").Value = 1
$ws.Range("B66").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C66").Value = "Arica y Parinacota"
$ws.Range("D66").Value = 44522
$ws.Range("D66").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E66").Value = 15
$ws.Range("F66").Value = 100112036
$ws.Range("G66").Value = "Caigua"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 120
$ws.Range("K66").Value = 4000
$ws.Range("L66").Value = 4500
$ws.Range("M66").Value = 4250
$ws.Range("N66").Value = "$/caja 20 kilos"
$ws.Range("O66").Value = "Región de Arica y Parinacota"
$ws.Range("P66").Value = 212
$ws.Range("Q66").Value = 20
$ws.Range("R66").Value = "Hortaliza"
